$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: explicit counts replaced with zeros (C26: 29 -> 0, D26: 175 -> 0).
# Downstream formulas (I26=SUM(C26:H26), etc.) recalc automatically.
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

# K33 / K35 change from hard-coded totals to formulas summing the row's
# four group columns. With all inputs still at 0 this evaluates to 0,
# which in turn makes the dependent percentage formulas in rows 34/36
# divide by zero (#DIV/0!) - that happens automatically on recalculation.
$ws.Range("K33").Formula = "=C33+E33+G33+I33"
$ws.Range("K35").Formula = "=C35+E35+G35+I35"

# Row 32 loses its explicit custom row height (back to sheet default).
$ws.Rows.Item(32).AutoFit()

# Last user selection moves from B6:K6 to the single cell D27.
$ws.Range("D27").Select() | Out-Null
